$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'baseball pants high knee'
$ws.Cells.Item(2, 1).Value = 'knee pack'
$ws.Cells.Item(3, 1).Value = 'girls compression pants'
$ws.Cells.Item(4, 1).Value = 'basketball gear for boys'
$ws.Cells.Item(5, 1).Value = 'gym knee pads'
$ws.Cells.Item(6, 1).Value = 'extra large knee pads'
$ws.Cells.Item(7, 1).Value = 'tall leggings'
$ws.Cells.Item(8, 1).Value = 'large knee pads'
$ws.Cells.Item(9, 1).Value = 'football shock pads'
$ws.Cells.Item(10, 1).Value = 'fitness compression pants men'
$ws.Cells.Item(11, 1).Value = 'spandex leggings'
$ws.Cells.Item(12, 1).Value = 'knee pad hiking'
$ws.Cells.Item(13, 1).Value = 'workout legging for men'
$ws.Cells.Item(14, 1).Value = 'compression legs men'
$ws.Cells.Item(15, 1).Value = 'pants for men'
$ws.Cells.Item(16, 1).Value = 'mens baseball clothing'
$ws.Cells.Item(17, 1).Value = 'capri legging pants'
$ws.Cells.Item(18, 1).Value = 'compression pants knee length'
$ws.Cells.Item(19, 1).Value = 'knee pad for pain'
$ws.Cells.Item(20, 1).Value = 'camo knee pads for basketball'
$ws.Cells.Item(21, 1).Value = 'mcdavid basketball tights knee pads'
$ws.Cells.Item(22, 1).Value = 'skateboard knee pads youth'
$ws.Cells.Item(23, 1).Value = 'kids basketball knee pads'
$ws.Cells.Item(24, 1).Value = 'youth basketball pants tearaway'
$ws.Cells.Item(25, 1).Value = 'kneee pads for basketball'
$ws.Cells.Item(26, 1).Value = 'commpression pants for basketball'
$ws.Cells.Item(27, 1).Value = 'supreme basketball leggings'
$ws.Cells.Item(28, 1).Value = 'volleyball knee pads youth girls'
$ws.Cells.Item(29, 1).Value = 'adidas knee pads basketball'
$ws.Cells.Item(30, 1).Value = 'basketball tights men'
$ws.Cells.Item(31, 1).Value = 'compression pants men adidas'
$ws.Cells.Item(32, 1).Value = 'pant with knee pads'
$ws.Cells.Item(33, 1).Value = 'knee pads tights'
$ws.Cells.Item(34, 1).Value = 'compression pants with knee'
$ws.Cells.Item(35, 1).Value = 'women basketball pants'
$ws.Cells.Item(36, 1).Value = 'basketball knee pads for youth kids'
$ws.Cells.Item(37, 1).Value = 'pantalones con rodilleras'
$ws.Cells.Item(38, 1).Value = 'caterpillar knee pad pants'
$ws.Cells.Item(39, 1).Value = 'thick mens leggings'
$ws.Cells.Item(40, 1).Value = 'capri leggings men'
$ws.Cells.Item(41, 1).Value = 'compression leggings men basketball'
$ws.Cells.Item(42, 1).Value = 'capri tights men'
$ws.Cells.Item(43, 1).Value = 'mens compression knee pads'
$ws.Cells.Item(44, 1).Value = 'legging men'
$ws.Cells.Item(45, 1).Value = 'basketball knee pads youth'
$ws.Cells.Item(46, 1).Value = 'hex pads basketball knee'
$ws.Cells.Item(47, 1).Value = 'baseball knee high pants'
$ws.Cells.Item(48, 1).Value = 'little boys baseball pants'
$ws.Cells.Item(49, 1).Value = 'boys baseball pants'
$ws.Cells.Item(50, 1).Value = 'knee pad compression'
$ws.Cells.Item(51, 1).Value = 'black compression tights'
$ws.Cells.Item(52, 1).Value = 'capris men'
$ws.Cells.Item(53, 1).Value = 'hip pads for volleyball'
$ws.Cells.Item(54, 1).Value = 'soccer pants men'
$ws.Cells.Item(55, 1).Value = 'boys knee pads'
$ws.Cells.Item(56, 1).Value = 'capris pants men'
$ws.Cells.Item(57, 1).Value = 'knee pad baseball'
$ws.Cells.Item(58, 1).Value = 'volleyball kneepads youth'
$ws.Cells.Item(59, 1).Value = 'compression leggings for boys'
$ws.Cells.Item(60, 1).Value = 'knee pads lightweight'
$ws.Cells.Item(61, 1).Value = 'knee pads sliding'
$ws.Cells.Item(62, 1).Value = 'capri leggings mesh'
$ws.Cells.Item(63, 1).Value = 'men running pants'
$ws.Cells.Item(64, 1).Value = 'knee pad black'
$ws.Cells.Item(65, 1).Value = 'baseball pants youth boys'
$ws.Cells.Item(66, 1).Value = 'patella protector'
$ws.Cells.Item(67, 1).Value = 'leggings pack'
$ws.Cells.Item(68, 1).Value = 'knee pads for workout'
$ws.Cells.Item(69, 1).Value = 'capri pants boys'
$ws.Cells.Item(70, 1).Value = 'gym pad men'
$ws.Cells.Item(71, 1).Value = 'mens athletic pants'
$ws.Cells.Item(72, 1).Value = 'leggings youth'
$ws.Cells.Item(73, 1).Value = 'knee pads protection'
$ws.Cells.Item(74, 1).Value = 'knee pad softball'
$ws.Cells.Item(75, 1).Value = 'knee pads fitness'
$ws.Cells.Item(76, 1).Value = 'leggings youth boys'
$ws.Cells.Item(77, 1).Value = 'gel knee pad'
$ws.Cells.Item(78, 1).Value = 'basketballs for boys'
$ws.Cells.Item(79, 1).Value = 'women athletic leggings'
$ws.Cells.Item(80, 1).Value = 'knee pads insert'
$ws.Cells.Item(81, 1).Value = 'capri leggings'
$ws.Cells.Item(82, 1).Value = 'elbow pads basketball youth'
$ws.Cells.Item(83, 1).Value = 'nike compression tights'
$ws.Cells.Item(84, 1).Value = 'nike compression leggings men'
$ws.Cells.Item(85, 1).Value = 'gloves and knee pads'
$ws.Cells.Item(86, 1).Value = 'black leggings running'
$ws.Cells.Item(87, 1).Value = 'lupo compression leggings'
$ws.Cells.Item(88, 1).Value = 'compression leggings 30'
$ws.Cells.Item(89, 1).Value = 'compression leggings black'
$ws.Cells.Item(90, 1).Value = 'compression leggings circulation'
$ws.Cells.Item(91, 1).Value = 'compression leggings girls'
$ws.Cells.Item(92, 1).Value = 'compression leggings men under armour'
$ws.Cells.Item(93, 1).Value = 'compression leggings men white'
$ws.Cells.Item(94, 1).Value = 'compression leggings men nike'
$ws.Cells.Item(95, 1).Value = 'compression leggings running'
$ws.Cells.Item(96, 1).Value = 'compression leggings xxl'
$ws.Cells.Item(97, 1).Value = 'running capri pants'
$ws.Cells.Item(98, 1).Value = 'knee basketball pads'
$ws.Cells.Item(99, 1).Value = 'girls black leggings'
$ws.Cells.Item(100, 1).Value = 'basketball knee pads kids youth'
